$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: Apply cell styles (format-only) for new columns BA:CU across rows 3-5 ---
$styleRef = @{ 10 = "C3"; 11 = "D3"; 12 = "K3" }
$colStyle = @{
    "BA" = 12;
    "BB" = 12;
    "BC" = 11;
    "BD" = 11;
    "BE" = 11;
    "BF" = 11;
    "BG" = 11;
    "BH" = 10;
    "BI" = 11;
    "BJ" = 11;
    "BK" = 11;
    "BL" = 11;
    "BM" = 11;
    "BN" = 11;
    "BO" = 10;
    "BP" = 11;
    "BQ" = 10;
    "BR" = 11;
    "BS" = 11;
    "BT" = 11;
    "BU" = 11;
    "BV" = 10;
    "BW" = 11;
    "BX" = 10;
    "BY" = 10;
    "BZ" = 10;
    "CA" = 10;
    "CB" = 10;
    "CC" = 10;
    "CD" = 12;
    "CE" = 11;
    "CF" = 12;
    "CG" = 12;
    "CH" = 12;
    "CI" = 12;
    "CJ" = 12;
    "CK" = 11;
    "CL" = 11;
    "CM" = 10;
    "CN" = 11;
    "CO" = 11;
    "CP" = 11;
    "CQ" = 11;
    "CR" = 11;
    "CS" = 11;
    "CT" = 11;
    "CU" = 10
}
foreach ($col in $colStyle.Keys) {
    $s = $colStyle[$col]
    $ref = $styleRef[$s]
    $ws.Range($ref).Copy()
    $ws.Range("$($col)3:$($col)5").PasteSpecial(-4122)
}

# --- Step 2: Set cell values ---
# Pass A: columns BA..CM, rows 3-5
# Row 3
$ws.Range("BA3").Value = 1234567890
$ws.Range("BB3").Value = 'name1'
$ws.Range("BC3").Value = 'abc123'
$ws.Range("BD3").Value = 'a'
$ws.Range("BE3").Value = 'j'
$ws.Range("BF3").Value = 'zx'
$ws.Range("BG3").Value = 'jk'
$ws.Range("BH3").Value = 123456789
$ws.Range("BI3").Value = 1234
$ws.Range("BJ3").Value = 'A'
$ws.Range("BK3").Value = 1234
$ws.Range("BL3").Value = 5112023
$ws.Range("BM3").Value = 18032024
$ws.Range("BN3").Value = 29092022
$ws.Range("BO3").Value = 7082023
$ws.Range("BP3").Value = 5112023
$ws.Range("BQ3").Value = 'ABCDEFG12345kkkkkkkkkkkkkkkkkk' + [char]10 + 'kAAAAAAAAAAAAAAAAAAAAAAAAAAAAA' + [char]10 + 'AAAAAAAAAAAAAAAAAAAAAAAAAAAAAA' + [char]10 + 'AAAAAAAAaaaaaaaaaaaaaaaaaaaaaa'
$ws.Range("BR3").Value = 1234
$ws.Range("BS3").Value = 1234567890
$ws.Range("BT3").Value = 123456789
$ws.Range("BU3").Value = 123456789
$ws.Range("BV3").Value = 1234
$ws.Range("BW3").Value = 'AB'
$ws.Range("BX3").Value = 'AAAAAAAAAAAAAAAAAAAAAAAAAAAAAA' + [char]10 + 'AAAAAAAAAAAAAAAAAAAAAAAAAAAAAA' + [char]10 + 'AAAAA'
$ws.Range("BY3").Value = 'AM'
$ws.Range("BZ3").Value = 123
$ws.Range("CA3").Value = 1234567890
$ws.Range("CB3").Value = 'nombre1'
$ws.Range("CC3").Value = 1234567890
$ws.Range("CD3").Value = 27042006
$ws.Range("CE3").Value = 1234567890
$ws.Range("CF3").Value = 1234567890
$ws.Range("CG3").Value = 20122009
$ws.Range("CH3").Value = 123456789
$ws.Range("CI3").Value = 98
$ws.Range("CJ3").Value = 123
$ws.Range("CK3").Value = 123
$ws.Range("CL3").Value = 1234567899
$ws.Range("CM3").Value = 'correo@gmail.com'
# Row 4
$ws.Range("BA4").Value = 123
$ws.Range("BB4").Value = 'name2'
$ws.Range("BC4").Value = 'cde456'
$ws.Range("BG4").Value = 'fr'
$ws.Range("BH4").Value = 1234
$ws.Range("BI4").Value = 123
$ws.Range("BJ4").Value = 'D'
$ws.Range("BK4").Value = 567
$ws.Range("BN4").Value = 22112022
$ws.Range("BO4").Value = 3072023
$ws.Range("BR4").Value = 5678
$ws.Range("BV4").Value = 9876
$ws.Range("BW4").Value = 'DE'
$ws.Range("BY4").Value = 'TK'
$ws.Range("BZ4").Value = 4
$ws.Range("CB4").Value = 'nombre2'
$ws.Range("CE4").Value = 123
$ws.Range("CF4").Value = 12345
$ws.Range("CG4").Value = 16012025
$ws.Range("CJ4").Value = 456
# Row 5
$ws.Range("BB5").Value = 'name3'
$ws.Range("BK5").Value = 8910
$ws.Range("BR5").Value = 123456789
$ws.Range("CB5").Value = 'nombre3'
$ws.Range("CJ5").Value = 8

# Pass B: columns CN..CU (only CT/CU are new here), rows 3-5
# Row 3
$ws.Range("CN3").Value = 8465498615
$ws.Range("CO3").Value = 12345678912345
$ws.Range("CP3").Value = 123456789
$ws.Range("CQ3").Value = 987654321
$ws.Range("CR3").Value = 123456789
$ws.Range("CS3").Value = 123456789
$ws.Range("CT3").Value = 'Arturo'
$ws.Range("CU3").Value = 'Mi domicilio muy largo para' + [char]10 + 'comprobar el salto de linea'
# Row 4
$ws.Range("CN4").Value = 65452543120651
$ws.Range("CT4").Value = 'Arturito'
# Row 5

# --- Step 3: Column widths (closest achievable given engine rounding) ---
$ws.Columns.Item(69).ColumnWidth = 38.6   # BQ -> target xml width 39.5703125 (closest achievable: 39.5)
$ws.Columns.Item(76).ColumnWidth = 38.6   # BX -> target xml width 39.5703125 (closest achievable: 39.5)
$ws.Columns.Item(99).ColumnWidth = 24.6   # CU -> target xml width 25.42578125 (closest achievable: 25.5)
